# Reflect the new "Discounted Total" row and refresh column A's width.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New label + formula row (row 7) beneath the existing "Total" row.
$ws.Range("A7").Value = "Discounted Total"
$ws.Range("B7").Formula = "= 90% * B6"

# Column A was best-fit/auto-sized to accommodate the new content.
$ws.Columns(1).ColumnWidth = 14

# Move the active selection to B8, just past the newly added row.
$ws.Range("B8").Select() | Out-Null
